# Add "Chris George-O'Neil (wow, this was complicated)" to the class list,
# and split "Renato G. Nascimento" into two runs around a spell-check
# proofErr marker on the surname, matching the captured Word edit.

$d = $word.ActiveDocument

# Smart/curly right single quotation mark used in "O'Neil" (Word
# autocorrects a typed straight apostrophe to this).
$rsquo = [char]0x2019

# Locate the "Renato G. Nascimento" paragraph (4th paragraph in the doc).
# Paragraph.Range.Text includes the trailing paragraph mark (chr 13), so
# trim it before comparing.
$target = $null
foreach ($p in $d.Paragraphs) {
    $ptext = $p.Range.Text.TrimEnd([char]13)
    if ($ptext -eq "Renato G. Nascimento") {
        $target = $p
    }
}
if ($target -eq $null) {
    # Fallback: known position in the original class-list template.
    $target = $d.Paragraphs.Item(4)
}

$openXmlNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$bodyFragment = "<w:p><w:r><w:t xml:space='preserve'>Renato G. </w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:t>Nascimento</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/></w:p>" +
    "<w:p><w:r><w:t>Chris George-O" + $rsquo + "Neil (wow, this was complicated)</w:t></w:r>" +
    "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>"

$packageXml = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" +
    "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
    "<pkg:xmlData><w:document $openXmlNs><w:body>$bodyFragment</w:body></w:document></pkg:xmlData>" +
    "</pkg:part></pkg:package>"

# Replace the whole "Renato G. Nascimento" paragraph (the bookmark moves
# along with the freshly typed name) with the split-run version plus the
# newly added paragraph for Chris George-O'Neil.
[void]$target.Range.InsertXML($packageXml)

# The document used to end with a blank placeholder paragraph; that is
# where "Chris George-O'Neil ..." was actually typed, so the now-redundant
# blank paragraph break right after it collapses away, leaving the
# trailing paragraph holding the new name (and the relocated bookmark).
$chrisPara = $d.Paragraphs.Item($target.Index + 1)
$mergeStart = $chrisPara.Range.End - 1
$d.Range($mergeStart, $mergeStart + 1).Delete()
